$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.806.57"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "2.569.67"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.35"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.78"
$ws.Range("E6").Value = "  -3.34%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.544"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.09"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.66"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("E13").Value = "  +6.50%  "
$ws.Range("D14").Value = "2.590.80"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.880"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.23"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "42.897.82"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "0.0₃0989"
$ws.Range("E18").Value = "  +1.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.63"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "253.27"
$ws.Range("E22").Value = "  -4.68%  "
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "28.84"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.29"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.87"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  -4.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.00"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.48"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("E33").Value = "  -5.13%  "
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.38"
$ws.Range("E36").Value = "  +10.71%  "
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.40"
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("B40").Value = "ApeXProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.13"
$ws.Range("E40").Value = "  +32.39%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0311"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.41"
$ws.Range("E42").Value = "  -3.74%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.88"
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("D44").Value = "2.089.35"
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.28"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.77"
$ws.Range("E47").Value = "  -4.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.94"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.78"
$ws.Range("E49").Value = "  +9.87%  "
$ws.Range("D50").Value = "2.818.05"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.67"
$ws.Range("E51").Value = "  +0.98%  "
